$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give A1 the same header style (bold + border + alignment) already used by B1:F1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Write the new header row, shifted one column to the left, with the
# MODEL_CONDITION label renamed to MODELCONDITION
$ws.Range("A1").Value = "EL_Phylonet50"
$ws.Range("B1").Value = "FNRATE_ASTRAL"
$ws.Range("C1").Value = "TAXON"
$ws.Range("D1").Value = "MODELCONDITION"
$ws.Range("E1").Value = "GENE"

# Write the data rows, shifted one column to the left
$ws.Range("A2").Value = 187
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "11-texon"
$ws.Range("D2").Value = "estimated_50genes_weakILS"
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = 187
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "11-texon"
$ws.Range("D3").Value = "estimated_50genes_weakILS"
$ws.Range("E3").Value = 2

# The old A2/A3 cells carried the header-ish border style; the values that
# now live there are plain data cells, so drop that formatting
$ws.Range("A2:A3").ClearFormats()

# Drop the now-unused sixth column so the sheet shrinks back to A:E
$ws.Range("F1:F3").Clear()
